$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.781.79"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.475.12"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'319.06"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "'93.05"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'0.0870"
$ws.Range("E10").Value = "  +9.84%  "
$ws.Range("D11").Value = "'33.37"
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "2.855.38"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "'6.92"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "'15.69"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "2.479.27"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "'0.797"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").Value = "41.708.59"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'71.12"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'11.33"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "'240.73"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").Value = "'2.76"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "'1.95"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'24.76"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "'36.57"
$ws.Range("E30").Value = "  +3.18%  "
$ws.Range("D31").Value = "'158.01"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "'5.50"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0765"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.58"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'17.49"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'1.87"
$ws.Range("E37").Value = "  +4.65%  "
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  +8.42%  "
$ws.Range("D42").Value = "'4.01"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "1.999.49"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").Value = "'18.98"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").Value = "'9.56"
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("D48").Value = "2.712.08"
$ws.Range("D49").Value = "'98.30"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "'74.93"
$ws.Range("E50").Value = "  +4.65%  "
$ws.Range("D51").Value = "'67.25"
$ws.Range("E51").Value = "  +0.08%  "